# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# The workbook tracked quarterly stock-holding snapshots for 301263-泰恩康.
# This adds a new "2022-Q4" worksheet (inserted before the existing
# "2022-Q3" tab) and records its summary line on the "总计" (totals) sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert a new worksheet named "2022-Q4" right before "2022-Q3"
# so the final tab order is: 总计, 2022-Q4, 2022-Q3, 2022-Q2
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newWs = $wb.Worksheets.Add($q3)
$newWs.Name = "2022-Q4"

# Match the look & feel (fonts/borders/alignment) of the other quarter
# sheets by copying their formats onto the new sheet.
$q3.Range("A1:H9").Copy()
$newWs.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Step 2: header row
# ---------------------------------------------------------------------------
$headerArr = New-Object 'object[,]' 1,7
$headerArr[0,0] = "基金代码"
$headerArr[0,1] = "基金名称"
$headerArr[0,2] = "基金规模"
$headerArr[0,3] = "股票总仓位"
$headerArr[0,4] = "仓位占比"
$headerArr[0,5] = "持有市值(亿元)"
$headerArr[0,6] = "仓位排名"
$newWs.Range("B1:H1").Value = $headerArr

# ---------------------------------------------------------------------------
# Step 3: fund holding rows (basic fund info as of 2022-Q4)
# columns: idx, code, name, scale, stockPosition, positionRatio, marketValue, rank
# marketValue of $null means the source recorded a literal numeric 0 instead
# of a "0.00" text value (rows 24/25 in the sheet).
# ---------------------------------------------------------------------------
$dataRows = @(
  ,@(0, "213008", "宝盈资源优选混合", "8.72", "89.15", "5.50", "0.4796", 6)
  ,@(1, "213006", "宝盈核心优势灵活配置混合A", "8.78", "77.93", "4.56", "0.4004", 6)
  ,@(2, "050201", "博时价值增长贰号混合", "9.54", "74.96", "2.52", "0.2404", 10)
  ,@(3, "011460", "鹏华创新成长混合A", "10.08", "78.91", "2.17", "0.2187", 10)
  ,@(4, "519170", "浦银安盛增长动力灵活配置混合A", "6.85", "85.43", "2.64", "0.1808", 5)
  ,@(5, "519110", "浦银安盛价值成长混合A", "6.74", "90.48", "2.65", "0.1786", 5)
  ,@(6, "016021", "华安优嘉精选混合A", "6.09", "61.72", "2.08", "0.1267", 7)
  ,@(7, "010383", "宝盈基础产业混合A", "2.81", "91.40", "4.48", "0.1259", 9)
  ,@(8, "016022", "华安优嘉精选混合C", "5.26", "61.72", "2.08", "0.1094", 7)
  ,@(9, "006377", "广发趋势动力灵活配置混合", "2.57", "87.77", "2.47", "0.0635", 10)
  ,@(10, "519113", "浦银安盛精致生活混合", "1.81", "90.69", "3.28", "0.0594", 3)
  ,@(11, "519120", "浦银安盛新兴产业混合A", "1.72", "91.65", "3.40", "0.0585", 3)
  ,@(12, "011471", "鹏华致远成长混合A", "1.67", "60.84", "1.86", "0.0311", 5)
  ,@(13, "002137", "诺安利鑫灵活配置混合A", "0.44", "89.87", "3.52", "0.0155", 8)
  ,@(14, "010384", "宝盈基础产业混合C", "0.21", "91.40", "4.48", "0.0094", 9)
  ,@(15, "000241", "宝盈核心优势灵活配置混合C", "0.20", "77.93", "4.56", "0.0091", 6)
  ,@(16, "011461", "鹏华创新成长混合C", "0.35", "78.91", "2.17", "0.0076", 10)
  ,@(17, "011472", "鹏华致远成长混合C", "0.06", "60.84", "1.86", "0.0011", 5)
  ,@(18, "014003", "浦银安盛增长动力灵活配置混合C", "0.03", "85.43", "2.64", "0.0008", 5)
  ,@(19, "014521", "诺安利鑫灵活配置混合C", "0.01", "89.87", "3.52", "0.0004", 8)
  ,@(20, "014011", "浦银安盛价值成长混合C", "0.01", "90.48", "2.65", "0.0003", 5)
  ,@(21, "002000", "工银新生利混合", "0.55", "28.98", "0.04", "0.0002", 10)
  ,@(22, "014061", "浦银安盛新兴产业混合C", "0.00", "91.65", "3.40", $null, 3)
  ,@(23, "960031", "浦银安盛价值成长混合H", "0.00", "90.48", "2.65", $null, 5)
)

$lastRow = 1 + $dataRows.Count   # 25

# Columns B..G hold text (fund code / name / percentages stored as strings,
# same as the other quarter sheets) even though most look numeric, so force
# text formatting before writing, then clear the format back to General
# (the cells keep plain/default style afterwards, matching the sibling
# sheets which never apply a custom number format here).
$textRange = $newWs.Range("B2:G" + $lastRow)
$textRange.NumberFormat = "@"

$r = 2
foreach ($row in $dataRows) {
    $newWs.Cells.Item($r, 1).Value = $row[0]
    $newWs.Cells.Item($r, 2).Value = $row[1]
    $newWs.Cells.Item($r, 3).Value = $row[2]
    $newWs.Cells.Item($r, 4).Value = $row[3]
    $newWs.Cells.Item($r, 5).Value = $row[4]
    $newWs.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq $null) {
        $newWs.Cells.Item($r, 7).NumberFormat = "General"
        $newWs.Cells.Item($r, 7).Style = "Normal"
        $newWs.Cells.Item($r, 7).Value = 0
    } else {
        $newWs.Cells.Item($r, 7).Value = $row[6]
    }
    $newWs.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$textRange.NumberFormat = "General"
$textRange.Style = "Normal"

# Column A keeps the same bold/centred index style used on the sibling
# quarter sheets (already copied above); just re-apply to be safe.
$idxRange = $newWs.Range("A2:A" + $lastRow)
$q3.Range("A2").Copy()
$idxRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Step 4: update the "总计" (totals) sheet — insert the new 2022-Q4 summary
# row at the top of the data, pushing 2022-Q3 / 2022-Q2 down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the new last row (A4) the same index-column style as the others
# before writing into it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(4, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 24
$total.Cells.Item(2, 4).Value = 2.32

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 8
$total.Cells.Item(3, 4).Value = 0.41

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 6
$total.Cells.Item(4, 4).Value = 0.09

Write-Output "2022-Q4 sheet added and totals updated"
